$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 151.875
$ws.Range("J6").Value = 219.33333
$ws.Range("L6").Value = 657.99999
$ws.Range("N6").Value = -881.99999

$ws.Range("H8").Value = 276.82224
$ws.Range("I8").Value = 49.5
$ws.Range("K8").Value = 148.5
$ws.Range("M8").Value = -9.5

$ws.Range("H70").Value = 15101897
$ws.Range("I70").Value = 10875717
$ws.Range("J70").Value = 19328076
$ws.Range("K70").Value = 32627151
$ws.Range("L70").Value = 57984228
$ws.Range("M70").Value = -32626881
$ws.Range("N70").Value = -57984768

$ws.Range("H73").Value = 15101897
$ws.Range("I73").Value = 10875717
$ws.Range("J73").Value = 19328076
$ws.Range("K73").Value = 32627151
$ws.Range("L73").Value = 57984228
$ws.Range("M73").Value = -32626215
$ws.Range("N73").Value = -57986100

$ws.Range("H88").Value = 15911545
$ws.Range("I88").Value = 41667100
$ws.Range("J88").Value = 61972.23
$ws.Range("K88").Value = 41667100
$ws.Range("L88").Value = 61972.23
$ws.Range("M88").Value = -41666694
$ws.Range("N88").Value = -62784.23

$ws.Range("H91").Value = 15911545
$ws.Range("I91").Value = 41667100
$ws.Range("J91").Value = 61972.23
$ws.Range("K91").Value = 41667100
$ws.Range("L91").Value = 61972.23
$ws.Range("M91").Value = -41665696
$ws.Range("N91").Value = -64780.23

$ws.Range("H116").Value = 41672332
$ws.Range("I116").Value = 41672332
$ws.Range("K116").Value = 41672332
$ws.Range("M116").Value = -41668890

$ws.Range("H125").Value = 83334210
$ws.Range("I125").Value = 125000140
$ws.Range("J125").Value = 2349.5
$ws.Range("K125").Value = 1125001260
$ws.Range("L125").Value = 21145.5
$ws.Range("M125").Value = -1124998800
$ws.Range("N125").Value = -26065.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 71430690
$ws.Range("I2").Value = 1492.6666
$ws.Range("K2").Value = 1492.6666
$ws.Range("M2").Value = -1379.6666

$ws.Range("H45").Value = 2178.7273
$ws.Range("I45").Value = 2253
$ws.Range("K45").Value = 2253
$ws.Range("M45").Value = -1876

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").Value = ""

$ws.Range("H74").Value = 29443.783
$ws.Range("I74").Value = 37085.465
$ws.Range("K74").Value = 37085.465
$ws.Range("M74").Value = -36211.465

$ws.Range("H77").Value = 29443.783
$ws.Range("I77").Value = 37085.465
$ws.Range("K77").Value = 185427.325
$ws.Range("M77").Value = -181059.325

$ws.Range("H97").Value = 3206872.5
$ws.Range("I97").Value = 1561.75
$ws.Range("K97").Value = 1561.75
$ws.Range("M97").Value = -1065.75

$ws.Range("H116").Value = 71430690
$ws.Range("I116").Value = 1492.6666
$ws.Range("K116").Value = 1492.6666
$ws.Range("M116").Value = 801.3334

$ws.Range("H132").Value = 5027.4746
$ws.Range("I132").Value = 3109.932
$ws.Range("J132").Value = 10652.267
$ws.Range("K132").Value = 9329.795999999998
$ws.Range("L132").Value = 31956.801
$ws.Range("M132").Value = -6799.795999999998
$ws.Range("N132").Value = -37016.801

$ws.Range("H139").Value = 66857.5
$ws.Range("J139").Value = 66857.5
$ws.Range("L139").Value = 66857.5
$ws.Range("N139").Value = -77137.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 71430690
$ws.Range("I3").Value = 1492.6666
$ws.Range("K3").Value = 1492.6666
$ws.Range("M3").Value = -1378.6666

$ws.Range("H20").Value = 11908202
$ws.Range("I20").Value = 20836342
$ws.Range("J20").Value = 4017.1667
$ws.Range("K20").Value = 20836342
$ws.Range("L20").Value = 4017.1667
$ws.Range("M20").Value = -20836095
$ws.Range("N20").Value = -4511.1667

$ws.Range("H94").Value = 949.1667
$ws.Range("I94").Value = 670.84375
$ws.Range("K94").Value = 670.84375
$ws.Range("M94").Value = -219.84375

$ws.Range("H99").Value = 8266108.5
$ws.Range("I99").Value = 1112
$ws.Range("J99").Value = 18184104
$ws.Range("K99").Value = 1112
$ws.Range("L99").Value = 18184104
$ws.Range("M99").Value = 386
$ws.Range("N99").Value = -18187100

$ws.Range("H105").Value = 1873.4
$ws.Range("I105").Value = 1634.75
$ws.Range("K105").Value = 1634.75
$ws.Range("M105").Value = 112.25

$ws.Range("H134").Value = 5958144.5
$ws.Range("I134").Value = 12501950
$ws.Range("J134").Value = 9230.182000000001
$ws.Range("K134").Value = 37505850
$ws.Range("L134").Value = 27690.546
$ws.Range("M134").Value = -37503315
$ws.Range("N134").Value = -32760.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 107.95
$ws.Range("I7").Value = 63.642857
$ws.Range("K7").Value = 63.642857
$ws.Range("M7").Value = 49.357143

$ws.Range("H22").Value = 611.36365
$ws.Range("J22").Value = 575
$ws.Range("L22").Value = 575
$ws.Range("N22").Value = -1275

$ws.Range("H31").Value = 6294
$ws.Range("I31").Value = 1700.1052
$ws.Range("J31").Value = 9203.467000000001
$ws.Range("K31").Value = 1700.1052
$ws.Range("L31").Value = 9203.467000000001
$ws.Range("M31").Value = -1405.1052
$ws.Range("N31").Value = -9793.467000000001

$ws.Range("H34").Value = 6294
$ws.Range("I34").Value = 1700.1052
$ws.Range("J34").Value = 9203.467000000001
$ws.Range("K34").Value = 1700.1052
$ws.Range("L34").Value = 9203.467000000001
$ws.Range("M34").Value = -1498.1052
$ws.Range("N34").Value = -9607.467000000001

$ws.Range("H132").Value = 3310.762
$ws.Range("I132").Value = 1794.9722
$ws.Range("J132").Value = 5331.815
$ws.Range("K132").Value = 5384.9166
$ws.Range("L132").Value = 15995.445
$ws.Range("M132").Value = -2854.9166
$ws.Range("N132").Value = -21055.445

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 12777.571
$ws.Range("I133").Value = 9888.6
$ws.Range("K133").Value = 29665.8
$ws.Range("M133").Value = -24605.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 86806.25
$ws.Range("J64").Value = 99075
$ws.Range("L64").Value = 99075
$ws.Range("N64").Value = -99571

$ws.Range("H67").Value = 86806.25
$ws.Range("J67").Value = 99075
$ws.Range("L67").Value = 99075
$ws.Range("N67").Value = -100791

$ws.Range("H132").Value = 2851.204
$ws.Range("I132").Value = 1523.5476
$ws.Range("J132").Value = 10817.143
$ws.Range("K132").Value = 4570.642800000001
$ws.Range("L132").Value = 32451.429
$ws.Range("M132").Value = -2040.642800000001
$ws.Range("N132").Value = -37511.429

$ws.Range("H139").Value = 66661.71000000001
$ws.Range("J139").Value = 66661.71000000001
$ws.Range("L139").Value = 66661.71000000001
$ws.Range("N139").Value = -76941.71000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3567.3635
$ws.Range("I100").Value = 2034.8
$ws.Range("K100").Value = 2034.8
$ws.Range("M100").Value = -1493.8

$ws.Range("H107").Value = 3109.4
$ws.Range("I107").Value = 3109.4
$ws.Range("K107").Value = 3109.4
$ws.Range("M107").Value = -1189.4

$ws.Range("H136").Value = 12756.19
$ws.Range("I136").Value = 2962.0588
$ws.Range("J136").Value = 19416.2
$ws.Range("K136").Value = 8886.1764
$ws.Range("L136").Value = 58248.60000000001
$ws.Range("M136").Value = -6336.1764
$ws.Range("N136").Value = -63348.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 28406.2
$ws.Range("J41").Value = 28686.5
$ws.Range("L41").Value = 28686.5
$ws.Range("N41").Value = -29466.5

$ws.Range("H96").Value = 2999
$ws.Range("I96").Value = 3003
$ws.Range("J96").Value = 2995
$ws.Range("K96").Value = 3003
$ws.Range("L96").Value = 2995
$ws.Range("M96").Value = -1630
$ws.Range("N96").Value = -5741

$ws.Range("H132").Value = 4180.241
$ws.Range("I132").Value = 5173.8184
$ws.Range("J132").Value = 1057.5714
$ws.Range("K132").Value = 15521.4552
$ws.Range("L132").Value = 3172.7142
$ws.Range("M132").Value = -12991.4552
$ws.Range("N132").Value = -8232.7142
